# Apply updated market-board price/profit figures to the Sargatanas_Profits sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 868.8889
$ws.Range("I19").Value = 293
$ws.Range("J19").Value = 1090.3846
$ws.Range("K19").Value = 293
$ws.Range("L19").Value = 1090.3846
$ws.Range("M19").Value = -118
$ws.Range("N19").Value = -1440.3846
$ws.Range("H40").Value = 4079.875
$ws.Range("I40").Value = 2474
$ws.Range("J40").Value = 5043.4
$ws.Range("K40").Value = 2474
$ws.Range("L40").Value = 5043.4
$ws.Range("M40").Value = -2299
$ws.Range("N40").Value = -5393.4
$ws.Range("H62").Value = 1100
$ws.Range("J62").Value = 1100
$ws.Range("L62").Value = 1100
$ws.Range("N62").Value = -2348
$ws.Range("H64").Value = 6368.769
$ws.Range("I64").Value = 6359.3228
$ws.Range("J64").Value = 6405.375
$ws.Range("K64").Value = 6359.3228
$ws.Range("L64").Value = 6405.375
$ws.Range("M64").Value = -6111.3228
$ws.Range("N64").Value = -6901.375
$ws.Range("H65").Value = 1100
$ws.Range("J65").Value = 1100
$ws.Range("L65").Value = 5500
$ws.Range("N65").Value = -11740
$ws.Range("H67").Value = 6368.769
$ws.Range("I67").Value = 6359.3228
$ws.Range("J67").Value = 6405.375
$ws.Range("K67").Value = 6359.3228
$ws.Range("L67").Value = 6405.375
$ws.Range("M67").Value = -5501.3228
$ws.Range("N67").Value = -8121.375
$ws.Range("H88").Value = 13923455
$ws.Range("I88").Value = 41668420
$ws.Range("J88").Value = 50973.125
$ws.Range("K88").Value = 41668420
$ws.Range("L88").Value = 50973.125
$ws.Range("M88").Value = -41668014
$ws.Range("N88").Value = -51785.125
$ws.Range("H91").Value = 13923455
$ws.Range("I91").Value = 41668420
$ws.Range("J91").Value = 50973.125
$ws.Range("K91").Value = 41668420
$ws.Range("L91").Value = 50973.125
$ws.Range("M91").Value = -41667016
$ws.Range("N91").Value = -53781.125
$ws.Range("H129").Value = 1532.9524
$ws.Range("I129").Value = 1007.9
$ws.Range("K129").Value = 3023.7
$ws.Range("M129").Value = 1976.3
$ws.Range("H132").Value = 1360.925
$ws.Range("I132").Value = 992.7143
$ws.Range("K132").Value = 2978.1429
$ws.Range("M132").Value = -448.1428999999998
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H137").Value = 5332.2285
$ws.Range("I137").Value = 3196.9412
$ws.Range("J137").Value = 7348.8887
$ws.Range("K137").Value = 9590.8236
$ws.Range("L137").Value = 22046.6661
$ws.Range("M137").Value = -7040.8236
$ws.Range("N137").Value = -27146.6661
$ws.Range("H138").Value = 1790400.5
$ws.Range("I138").Value = 2532.45
$ws.Range("J138").Value = 2783660.5
$ws.Range("K138").Value = 7597.349999999999
$ws.Range("L138").Value = 8350981.5
$ws.Range("M138").Value = -2457.349999999999
$ws.Range("N138").Value = -8361261.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2707990.5
$ws.Range("I32").Value = 2860947
$ws.Range("K32").Value = 2860947
$ws.Range("M32").Value = -2860660
$ws.Range("H61").Value = 55558960
$ws.Range("I61").Value = 1572.4546
$ws.Range("J61").Value = 142863420
$ws.Range("K61").Value = 1572.4546
$ws.Range("L61").Value = 142863420
$ws.Range("M61").Value = -1360.4546
$ws.Range("N61").Value = -142863844
$ws.Range("H107").Value = 57777.5
$ws.Range("J107").Value = 57777.5
$ws.Range("L107").Value = 57777.5
$ws.Range("N107").Value = -65457.5
$ws.Range("H110").Value = 18525766
$ws.Range("I110").Value = 9039.357
$ws.Range("K110").Value = 9039.357
$ws.Range("M110").Value = -6994.357
$ws.Range("H122").Value = 2672.8386
$ws.Range("I122").Value = 1945.4073
$ws.Range("K122").Value = 5836.2219
$ws.Range("M122").Value = -3386.2219
$ws.Range("H132").Value = 4144.577
$ws.Range("I132").Value = 2687.2104
$ws.Range("K132").Value = 8061.6312
$ws.Range("M132").Value = -5531.6312
$ws.Range("H136").Value = 55558960
$ws.Range("I136").Value = 1572.4546
$ws.Range("J136").Value = 142863420
$ws.Range("K136").Value = 4717.3638
$ws.Range("L136").Value = 428590260
$ws.Range("M136").Value = -2167.3638
$ws.Range("N136").Value = -428595360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 38022.81
$ws.Range("J95").Value = 38022.81
$ws.Range("L95").Value = 38022.81
$ws.Range("N95").Value = -43514.81
$ws.Range("H107").Value = 102283384
$ws.Range("I107").Value = 160728910
$ws.Range("K107").Value = 160728910
$ws.Range("M107").Value = -160726990
$ws.Range("H112").Value = 59343
$ws.Range("J112").Value = 59343
$ws.Range("L112").Value = 59343
$ws.Range("N112").Value = -62297
$ws.Range("H117").Value = 95900
$ws.Range("J117").Value = 95900
$ws.Range("L117").Value = 95900
$ws.Range("N117").Value = -105078
$ws.Range("H134").Value = 4240373.5
$ws.Range("I134").Value = 5210319.5
$ws.Range("K134").Value = 15630958.5
$ws.Range("M134").Value = -15628423.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 248.45454
$ws.Range("I7").Value = 97.2
$ws.Range("J7").Value = 374.5
$ws.Range("K7").Value = 97.2
$ws.Range("L7").Value = 374.5
$ws.Range("M7").Value = 15.8
$ws.Range("N7").Value = -600.5
$ws.Range("H16").Value = 5514.778
$ws.Range("I16").Value = 738.6
$ws.Range("K16").Value = 738.6
$ws.Range("M16").Value = -451.6
$ws.Range("H31").Value = 4162.7544
$ws.Range("I31").Value = 1693.5555
$ws.Range("J31").Value = 6385.033
$ws.Range("K31").Value = 1693.5555
$ws.Range("L31").Value = 6385.033
$ws.Range("M31").Value = -1398.5555
$ws.Range("N31").Value = -6975.033
$ws.Range("H34").Value = 4162.7544
$ws.Range("I34").Value = 1693.5555
$ws.Range("J34").Value = 6385.033
$ws.Range("K34").Value = 1693.5555
$ws.Range("L34").Value = 6385.033
$ws.Range("M34").Value = -1491.5555
$ws.Range("N34").Value = -6789.033
$ws.Range("H51").Value = 27641.455
$ws.Range("I51").Value = 5090
$ws.Range("J51").Value = 32652.889
$ws.Range("K51").Value = 5090
$ws.Range("L51").Value = 32652.889
$ws.Range("M51").Value = -4354
$ws.Range("N51").Value = -34124.889
$ws.Range("H58").Value = 4779.2646
$ws.Range("I58").Value = 2470.0625
$ws.Range("K58").Value = 2470.0625
$ws.Range("M58").Value = -2267.0625
$ws.Range("H61").Value = 27641.455
$ws.Range("I61").Value = 5090
$ws.Range("J61").Value = 32652.889
$ws.Range("K61").Value = 5090
$ws.Range("L61").Value = 32652.889
$ws.Range("M61").Value = -4742
$ws.Range("N61").Value = -33348.889
$ws.Range("H86").Value = 24044020
$ws.Range("I86").Value = 31256098
$ws.Range("J86").Value = 3766.6667
$ws.Range("K86").Value = 31256098
$ws.Range("L86").Value = 3766.6667
$ws.Range("M86").Value = -31254975
$ws.Range("N86").Value = -6012.6667
$ws.Range("H89").Value = 24044020
$ws.Range("I89").Value = 31256098
$ws.Range("J89").Value = 3766.6667
$ws.Range("K89").Value = 156280490
$ws.Range("L89").Value = 18833.3335
$ws.Range("M89").Value = -156274874
$ws.Range("N89").Value = -30065.3335
$ws.Range("H113").Value = 5514.778
$ws.Range("I113").Value = 738.6
$ws.Range("K113").Value = 738.6
$ws.Range("M113").Value = 1431.4
$ws.Range("H122").Value = 3080.5
$ws.Range("I122").Value = 1611.6
$ws.Range("K122").Value = 4834.799999999999
$ws.Range("M122").Value = -2384.799999999999
$ws.Range("H132").Value = 5477.2085
$ws.Range("I132").Value = 4213.4287
$ws.Range("K132").Value = 12640.2861
$ws.Range("M132").Value = -10110.2861
$ws.Range("H136").Value = 4779.2646
$ws.Range("I136").Value = 2470.0625
$ws.Range("K136").Value = 7410.1875
$ws.Range("M136").Value = -4860.1875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 86810.75999999999
$ws.Range("I2").Value = 14391.667
$ws.Range("K2").Value = 86350.00199999999
$ws.Range("M2").Value = -86237.00199999999
$ws.Range("H4").Value = 73920980
$ws.Range("J4").Value = 4337904.5
$ws.Range("L4").Value = 13013713.5
$ws.Range("N4").Value = -13013937.5
$ws.Range("H11").Value = 2134.6667
$ws.Range("I11").Value = 452
$ws.Range("J11").Value = 5500
$ws.Range("K11").Value = 1356
$ws.Range("L11").Value = 16500
$ws.Range("M11").Value = -1216
$ws.Range("N11").Value = -16780
$ws.Range("H12").Value = 2381404.2
$ws.Range("I12").Value = 1248.5714
$ws.Range("J12").Value = 3571482.2
$ws.Range("K12").Value = 3745.7142
$ws.Range("L12").Value = 10714446.6
$ws.Range("M12").Value = -3572.7142
$ws.Range("N12").Value = -10714792.6
$ws.Range("H23").Value = 176.36363
$ws.Range("I23").Value = 88.333336
$ws.Range("K23").Value = 265.000008
$ws.Range("M23").Value = -30.00000799999998
$ws.Range("H33").Value = 20833608
$ws.Range("I33").Value = 55555636
$ws.Range("J33").Value = 389.5
$ws.Range("K33").Value = 333333816
$ws.Range("L33").Value = 2337
$ws.Range("M33").Value = -333333533
$ws.Range("N33").Value = -2903
$ws.Range("H34").Value = 5236.3335
$ws.Range("J34").Value = 5494.55
$ws.Range("L34").Value = 16483.65
$ws.Range("N34").Value = -16651.65
$ws.Range("H38").Value = 36.714287
$ws.Range("I38").Value = 26.166666
$ws.Range("J38").Value = 100
$ws.Range("K38").Value = 78.49999800000001
$ws.Range("L38").Value = 300
$ws.Range("M38").Value = 268.500002
$ws.Range("N38").Value = -994
$ws.Range("H86").Value = 618
$ws.Range("I86").Value = 112.333336
$ws.Range("K86").Value = 337.000008
$ws.Range("M86").Value = 848.999992
$ws.Range("H89").Value = 618
$ws.Range("I89").Value = 112.333336
$ws.Range("K89").Value = 1011.000024
$ws.Range("M89").Value = 4916.999976
$ws.Range("H97").Value = 690.6667
$ws.Range("J97").Value = 474
$ws.Range("L97").Value = 1422
$ws.Range("N97").Value = -2414
$ws.Range("H113").Value = 775
$ws.Range("I113").Value = 683.3333
$ws.Range("J113").Value = 800
$ws.Range("K113").Value = 2049.9999
$ws.Range("L113").Value = 2400
$ws.Range("M113").Value = 120.0001000000002
$ws.Range("N113").Value = -6740
$ws.Range("H137").Value = 89716.69500000001
$ws.Range("J137").Value = 69486.664
$ws.Range("L137").Value = 208459.992
$ws.Range("N137").Value = -218659.992
$ws.Range("H141").Value = 7594.125
$ws.Range("I141").Value = 5745.222
$ws.Range("J141").Value = 9971.286
$ws.Range("K141").Value = 17235.666
$ws.Range("L141").Value = 29913.858
$ws.Range("M141").Value = -12055.666
$ws.Range("N141").Value = -40273.858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 3998
$ws.Range("J10").Value = 3998
$ws.Range("L10").Value = 3998
$ws.Range("N10").Value = -4336
$ws.Range("H31").Value = 2500
$ws.Range("I31").Value = 2500
$ws.Range("K31").Value = 2500
$ws.Range("M31").Value = -2208
$ws.Range("H37").Value = 2500
$ws.Range("I37").Value = 2500
$ws.Range("K37").Value = 2500
$ws.Range("M37").Value = -2223
$ws.Range("H70").Value = 8058.143
$ws.Range("I70").Value = 5435.5
$ws.Range("J70").Value = 10025.125
$ws.Range("K70").Value = 5435.5
$ws.Range("L70").Value = 10025.125
$ws.Range("M70").Value = -5165.5
$ws.Range("N70").Value = -10565.125
$ws.Range("H73").Value = 8058.143
$ws.Range("I73").Value = 5435.5
$ws.Range("J73").Value = 10025.125
$ws.Range("K73").Value = 5435.5
$ws.Range("L73").Value = 10025.125
$ws.Range("M73").Value = -4499.5
$ws.Range("N73").Value = -11897.125
$ws.Range("H92").Value = 19440
$ws.Range("J92").Value = 19440
$ws.Range("L92").Value = 19440
$ws.Range("N92").Value = -23184
$ws.Range("H102").Value = 4009.2363
$ws.Range("I102").Value = 3865.6875
$ws.Range("J102").Value = 4993.5713
$ws.Range("K102").Value = 3865.6875
$ws.Range("L102").Value = 4993.5713
$ws.Range("M102").Value = -2243.6875
$ws.Range("N102").Value = -8237.5713
$ws.Range("H113").Value = 6721.3516
$ws.Range("I113").Value = 4588.5557
$ws.Range("K113").Value = 4588.5557
$ws.Range("M113").Value = -2418.5557
$ws.Range("H122").Value = 3461993
$ws.Range("J122").Value = 3866.6667
$ws.Range("L122").Value = 11600.0001
$ws.Range("N122").Value = -16500.0001
$ws.Range("H132").Value = 2457.9656
$ws.Range("I132").Value = 2113.48
$ws.Range("K132").Value = 6340.440000000001
$ws.Range("M132").Value = -3810.440000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5582.4287
$ws.Range("I7").Value = 4128.3335
$ws.Range("K7").Value = 4128.3335
$ws.Range("M7").Value = -4016.3335
$ws.Range("H126").Value = 5582.4287
$ws.Range("I126").Value = 4128.3335
$ws.Range("K126").Value = 12385.0005
$ws.Range("M126").Value = -9915.000499999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 254980.12
$ws.Range("I122").Value = 447278.34
$ws.Range("J122").Value = 7739.5713
$ws.Range("K122").Value = 1341835.02
$ws.Range("L122").Value = 23218.7139
$ws.Range("M122").Value = -1339385.02
$ws.Range("N122").Value = -28118.7139
$ws.Range("H126").Value = 1000
$ws.Range("I126").Value = 1000
$ws.Range("K126").Value = 3000
$ws.Range("M126").Value = -530
$ws.Range("H132").Value = 7119.5713
$ws.Range("I132").Value = 6719.5625
$ws.Range("J132").Value = 8399.6
$ws.Range("K132").Value = 20158.6875
$ws.Range("L132").Value = 25198.8
$ws.Range("M132").Value = -17628.6875
$ws.Range("N132").Value = -30258.8
$ws.Range("H136").Value = 211648.42
$ws.Range("I136").Value = 648.6667
$ws.Range("J136").Value = 482933.8
$ws.Range("K136").Value = 1946.0001
$ws.Range("L136").Value = 1448801.4
$ws.Range("M136").Value = 603.9999
$ws.Range("N136").Value = -1453901.4
